$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the old last row (35) down by 5 rows, opening up rows 35-39 for
# the new weather readings (old row 35 becomes row 40, unchanged except
# for the two cells edited below).
$ws.Rows("35:39").Insert()

# Dates in column C must stay plain text (matching the rest of the sheet)
# rather than being auto-parsed into Excel date serials.
$ws.Range("C35:C39").NumberFormat = "@"

$ws.Range("A35").Value = 28
$ws.Range("B35").Value = "Partly Cloudy"
$ws.Range("C35").Value = "01/18/2025"
$ws.Range("D35").Value = 18

$ws.Range("A36").Value = 28
$ws.Range("B36").Value = "Partly Cloudy"
$ws.Range("C36").Value = "01/18/2025"
$ws.Range("D36").Value = 18

$ws.Range("A37").Value = 28
$ws.Range("B37").Value = "Partly Cloudy"
$ws.Range("C37").Value = "01/18/2025"
$ws.Range("D37").Value = 18

$ws.Range("A38").Value = 28
$ws.Range("B38").Value = "Partly Cloudy"
$ws.Range("C38").Value = "01/18/2025"
$ws.Range("D38").Value = 19

$ws.Range("A39").Value = 28
$ws.Range("B39").Value = "Partly Cloudy"
$ws.Range("C39").Value = "01/18/2025"
$ws.Range("D39").Value = 19

# Drop the temporary text format again so these cells end up with the
# same (default) style as the rest of the data rows.
$ws.Range("C35:C39").ClearFormats()

# Old row 35 (now row 40) keeps its Weather Condition/Date, but the
# Temperature and Current Hour readings change; keep them as text like
# the original row.
$ws.Range("A40").NumberFormat = "@"
$ws.Range("A40").Value = "27"
$ws.Range("A40").ClearFormats()

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19"
$ws.Range("D40").ClearFormats()
